$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "#100009"
$ws.Range("B11").Value = "CreateMain Page"

$ws.Range("B12").Select()
